# Updates crypto price/volume data to the latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to Text format so numeric-looking strings (and
    # percentages) are stored verbatim instead of being reinterpreted
    # as numbers, then drop back to the Normal style so no stray
    # number-format override is left on the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "300.59"
Set-TextValue $ws.Range("E2") "-1.13%"
Set-TextValue $ws.Range("D3") "31.50"
Set-TextValue $ws.Range("E3") "-1.45%"
Set-TextValue $ws.Range("D4") "5.147"
Set-TextValue $ws.Range("E4") "-3.02%"
Set-TextValue $ws.Range("D5") "0.07355"
Set-TextValue $ws.Range("E5") "-1.46%"
Set-TextValue $ws.Range("D6") "1.839"
Set-TextValue $ws.Range("E6") "26.75%"
Set-TextValue $ws.Range("D7") "7.790"
Set-TextValue $ws.Range("E7") "-0.23%"
Set-TextValue $ws.Range("D8") "3.749"
Set-TextValue $ws.Range("E8") "-1.53%"
Set-TextValue $ws.Range("D9") "0.9278"
Set-TextValue $ws.Range("E9") "0.83%"
Set-TextValue $ws.Range("D10") "0.1693"
Set-TextValue $ws.Range("E10") "-0.44%"
Set-TextValue $ws.Range("D11") "0.07019"
Set-TextValue $ws.Range("E11") "-8.42%"
Set-TextValue $ws.Range("D12") "0.08041"
Set-TextValue $ws.Range("E12") "-0.52%"
Set-TextValue $ws.Range("D13") "0.03028"
Set-TextValue $ws.Range("E13") "0.25%"
Set-TextValue $ws.Range("D14") "0.09936"
Set-TextValue $ws.Range("E14") "0.58%"
Set-TextValue $ws.Range("D15") "0.001485"
Set-TextValue $ws.Range("E15") "-0.34%"
Set-TextValue $ws.Range("D16") "0.006150"
Set-TextValue $ws.Range("E16") "-3.80%"
Set-TextValue $ws.Range("D17") "3.460"
Set-TextValue $ws.Range("E17") "-0.63%"
Set-TextValue $ws.Range("D18") "2.225"
Set-TextValue $ws.Range("E18") "-0.12%"
Set-TextValue $ws.Range("D19") "0.3274"
Set-TextValue $ws.Range("E19") "-1.78%"
Set-TextValue $ws.Range("D20") "0.1329"
Set-TextValue $ws.Range("E20") "-1.23%"
Set-TextValue $ws.Range("D21") "4.554"
Set-TextValue $ws.Range("E21") "1.63%"
Set-TextValue $ws.Range("D22") "0.04641"
Set-TextValue $ws.Range("D23") "0.1579"
Set-TextValue $ws.Range("E23") "-2.73%"
Set-TextValue $ws.Range("E24") "-0.38%"
Set-TextValue $ws.Range("D25") "0.004759"
Set-TextValue $ws.Range("E25") "7.93%"
Set-TextValue $ws.Range("D26") "0.0001297"
Set-TextValue $ws.Range("E26") "-7.37%"
Set-TextValue $ws.Range("E27") "7.38%"
Set-TextValue $ws.Range("D39") "0.01720"
Set-TextValue $ws.Range("E39") "0.03%"
Set-TextValue $ws.Range("D40") "0.04491"
Set-TextValue $ws.Range("E40") "-0.62%"
Set-TextValue $ws.Range("D41") "0.007090"
Set-TextValue $ws.Range("E41") "-1.59%"
Set-TextValue $ws.Range("D42") "0.1340"
Set-TextValue $ws.Range("E42") "0.03%"
Set-TextValue $ws.Range("D43") "0.002166"
Set-TextValue $ws.Range("E43") "-3.36%"
Set-TextValue $ws.Range("D44") "0.01098"
Set-TextValue $ws.Range("E44") "-13.02%"
Set-TextValue $ws.Range("D45") "0.00006219"
Set-TextValue $ws.Range("E45") "1.14%"
Set-TextValue $ws.Range("E46") "-21.52%"
Set-TextValue $ws.Range("D47") "0.7397"
Set-TextValue $ws.Range("E47") "-60.50%"
